$wb = $excel.ActiveWorkbook

# --- ELF-bldg-winter sheet ---
$wsWinter = $wb.Worksheets.Item("ELF-bldg-winter")
$wsWinter.Range("D2").Value = 3.21241
$wsWinter.Range("B5").Value = 1.43078
$wsWinter.Range("D5").Value = 1.05923
$wsWinter.Range("D7").Value = 1.05923

# --- ELF-bldg-summer sheet ---
$wsSummer = $wb.Worksheets.Item("ELF-bldg-summer")
$wsSummer.Range("B3").Value = 2.76336
$wsSummer.Range("D3").Value = 4.78345
$wsSummer.Range("B5").Value = 1.46815
$wsSummer.Range("D5").Value = 1.694
$wsSummer.Range("D7").Value = 1.694

# --- ELF-vehicles sheet ---
$wsVehicles = $wb.Worksheets.Item("ELF-vehicles")
$wsVehicles.Range("B4").Value = 1.0282
$wsVehicles.Range("C4").Value = 0.7972399999999999
$wsVehicles.Range("B5").Value = 1.0282
$wsVehicles.Range("C5").Value = 0.7972399999999999
$wsVehicles.Range("B6").Value = 1.0282
$wsVehicles.Range("C6").Value = 0.7972399999999999
$wsVehicles.Range("B7").Value = 1.0282
$wsVehicles.Range("C7").Value = 0.7972399999999999
